$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 37
$ws.Range("I2").Value = 103
$ws.Range("J2").Value = 407
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 119
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 58
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 35
$ws.Range("T2").Value = 63
$ws.Range("U2").Value = 8
$ws.Range("V2").Value = 574
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 612
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 12
